$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rich text edits: header "Volume N  Number NN" and "Report Covering the Week ... Through ..." ---
# Only the specific substring runs change; other runs/formatting are left untouched.
$ws.Range("A8").Characters(21,2).Text = "44"
$ws.Range("C9").Characters(27,10).Text = "10/30/2023"
$ws.Range("C9").Characters(48,10).Text = "11/5/2023"

# --- Structural changes: cells that flip between a numeric literal and the
#     shared "0" / "***.*" text placeholders used elsewhere on this sheet.
#     We copy number-format only from a donor cell that already has the right
#     look, then push the value, so style/type both land correctly. ---

# Row 22 (Transit): C/D/E swap between numeric and text roles
$ws.Range("F22").Copy()
$ws.Range("C22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E22").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("E22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 28 (Shooting Vic.): D/E become text placeholders
$ws.Range("D28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E28").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("E28").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 29 (Shooting Inc.): D/E become text placeholders
$ws.Range("D29").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E29").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("E29").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 30 (Hate Crimes): D/E become text placeholders
$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E30").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("E30").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Plain numeric literal updates across the crime table ---
$ws.Range("L14").Value = 25
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = -42.857142857142
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -14.285714285714
$ws.Range("I16").Value = 178
$ws.Range("J16").Value = 164
$ws.Range("K16").Value = 8.536585365853
$ws.Range("L16").Value = -3.260869565217
$ws.Range("M16").Value = -52.150537634408
$ws.Range("N16").Value = -89.874857792946
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 37
$ws.Range("H17").Value = -32.432432432432
$ws.Range("I17").Value = 354
$ws.Range("J17").Value = 395
$ws.Range("K17").Value = -10.379746835443
$ws.Range("L17").Value = -1.392757660167
$ws.Range("M17").Value = -3.542234332425
$ws.Range("N17").Value = -64.242424242424
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 202
$ws.Range("J18").Value = 211
$ws.Range("K18").Value = -4.265402843601
$ws.Range("L18").Value = 3.061224489795
$ws.Range("M18").Value = -44.505494505494
$ws.Range("N18").Value = -81.209302325581
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = -60
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -18.60465116279
$ws.Range("I19").Value = 383
$ws.Range("J19").Value = 353
$ws.Range("K19").Value = 8.498583569405
$ws.Range("L19").Value = 15.709969788519
$ws.Range("M19").Value = 6.983240223463
$ws.Range("N19").Value = -48.243243243243
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 16
$ws.Range("H20").Value = 300
$ws.Range("I20").Value = 110
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = 34.146341463414
$ws.Range("L20").Value = 37.5
$ws.Range("M20").Value = 7.843137254901
$ws.Range("N20").Value = -76.344086021505
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 33
$ws.Range("E21").Value = -57.575757575757
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 118
$ws.Range("H21").Value = -9.322033898305
$ws.Range("I21").Value = 1257
$ws.Range("J21").Value = 1236
$ws.Range("K21").Value = 1.699029126213
$ws.Range("L21").Value = 7.344150298889
$ws.Range("M21").Value = -21.388367729831
$ws.Range("N21").Value = -75.549503987551
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = -75
$ws.Range("I22").Value = 19
$ws.Range("K22").Value = -20.833333333333
$ws.Range("L22").Value = -20.833333333333
$ws.Range("M22").Value = -36.666666666666
$ws.Range("F23").Value = 14
$ws.Range("H23").Value = -6.666666666666
$ws.Range("I23").Value = 211
$ws.Range("J23").Value = 205
$ws.Range("K23").Value = 2.926829268292
$ws.Range("L23").Value = -1.401869158878
$ws.Range("M23").Value = 12.234042553191
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = -30.76923076923
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 148
$ws.Range("H24").Value = -22.972972972973
$ws.Range("I24").Value = 1284
$ws.Range("J24").Value = 1392
$ws.Range("K24").Value = -7.758620689655
$ws.Range("L24").Value = 16.409791477787
$ws.Range("M24").Value = 49.302325581395
$ws.Range("C25").Value = 13
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = -43.478260869565
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = -18.333333333333
$ws.Range("I25").Value = 544
$ws.Range("J25").Value = 580
$ws.Range("K25").Value = -6.206896551724
$ws.Range("L25").Value = 25.345622119815
$ws.Range("M25").Value = -35.697399527186
$ws.Range("C27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 43
$ws.Range("J27").Value = 36
$ws.Range("K27").Value = 19.444444444444
$ws.Range("L27").Value = -38.571428571428
$ws.Range("L28").Value = -40.74074074074
$ws.Range("N28").Value = -86.991869918699
$ws.Range("L29").Value = -46.511627906976
$ws.Range("N29").Value = -89.302325581395
